$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 166.25
$ws.Range("I42").Value = 55
$ws.Range("K42").Value = 165
$ws.Range("M42").Value = 65
$ws.Range("H86").Value = 3028.0908
$ws.Range("I86").Value = 2749.6667
$ws.Range("J86").Value = 3624.7144
$ws.Range("K86").Value = 2749.6667
$ws.Range("L86").Value = 3624.7144
$ws.Range("M86").Value = -1626.6667
$ws.Range("N86").Value = -5870.7144
$ws.Range("H89").Value = 3028.0908
$ws.Range("I89").Value = 2749.6667
$ws.Range("J89").Value = 3624.7144
$ws.Range("K89").Value = 13748.3335
$ws.Range("L89").Value = 18123.572
$ws.Range("M89").Value = -8132.333500000001
$ws.Range("N89").Value = -29355.572
$ws.Range("H98").Value = 1515.8572
$ws.Range("I98").Value = 1586.4166
$ws.Range("J98").Value = 1092.5
$ws.Range("K98").Value = 1586.4166
$ws.Range("L98").Value = 1092.5
$ws.Range("M98").Value = -88.41660000000002
$ws.Range("N98").Value = -4088.5
$ws.Range("H122").Value = 1515.8572
$ws.Range("I122").Value = 1586.4166
$ws.Range("J122").Value = 1092.5
$ws.Range("K122").Value = 4759.2498
$ws.Range("L122").Value = 3277.5
$ws.Range("M122").Value = -2309.2498
$ws.Range("N122").Value = -8177.5
$ws.Range("H125").Value = 692
$ws.Range("I125").Value = 692
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 6228
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -3768
$ws.Range("N125").ClearContents()
$ws.Range("H138").Value = 2406.2222
$ws.Range("J138").Value = 2949.6316
$ws.Range("L138").Value = 8848.8948
$ws.Range("N138").Value = -19128.8948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 20681.9
$ws.Range("I74").Value = 22008.117
$ws.Range("K74").Value = 22008.117
$ws.Range("M74").Value = -21134.117
$ws.Range("H77").Value = 20681.9
$ws.Range("I77").Value = 22008.117
$ws.Range("K77").Value = 110040.585
$ws.Range("M77").Value = -105672.585
$ws.Range("H132").Value = 1612.3334
$ws.Range("I132").Value = 1412.1034
$ws.Range("K132").Value = 4236.3102
$ws.Range("M132").Value = -1706.3102

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3789.8333
$ws.Range("I31").Value = 4082.875
$ws.Range("J31").Value = 3643.3125
$ws.Range("K31").Value = 4082.875
$ws.Range("L31").Value = 3643.3125
$ws.Range("M31").Value = -3787.875
$ws.Range("N31").Value = -4233.3125
$ws.Range("H34").Value = 3789.8333
$ws.Range("I34").Value = 4082.875
$ws.Range("J34").Value = 3643.3125
$ws.Range("K34").Value = 4082.875
$ws.Range("L34").Value = 3643.3125
$ws.Range("M34").Value = -3880.875
$ws.Range("N34").Value = -4047.3125
$ws.Range("H107").Value = 1530.4445
$ws.Range("I107").Value = 1574.1428
$ws.Range("J107").Value = 1377.5
$ws.Range("K107").Value = 1574.1428
$ws.Range("L107").Value = 1377.5
$ws.Range("M107").Value = 345.8571999999999
$ws.Range("N107").Value = -5217.5
$ws.Range("H132").Value = 26485
$ws.Range("J132").Value = 32335.385
$ws.Range("L132").Value = 97006.155
$ws.Range("N132").Value = -102066.155
$ws.Range("H134").Value = 4202.722
$ws.Range("I134").Value = 2343.5334
$ws.Range("K134").Value = 7030.600199999999
$ws.Range("M134").Value = -4495.600199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2608.3489
$ws.Range("I131").Value = 1097.8572
$ws.Range("K131").Value = 3293.5716
$ws.Range("M131").Value = 1746.4284
$ws.Range("H137").Value = 5708.875
$ws.Range("I137").Value = 2947.6
$ws.Range("J137").Value = 10311
$ws.Range("K137").Value = 8842.799999999999
$ws.Range("L137").Value = 30933
$ws.Range("M137").Value = -3742.799999999999
$ws.Range("N137").Value = -41133
$ws.Range("H139").Value = 2537.6667
$ws.Range("I139").Value = 2537.6667
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 7613.000100000001
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -2473.000100000001
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 44890.1
$ws.Range("J46").Value = 44890.1
$ws.Range("L46").Value = 44890.1
$ws.Range("N46").Value = -45202.1
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H132").Value = 15281.972
$ws.Range("I132").Value = 14589.261
$ws.Range("K132").Value = 43767.783
$ws.Range("M132").Value = -41237.783

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 60647.6
$ws.Range("J6").Value = 60647.6
$ws.Range("L6").Value = 60647.6
$ws.Range("N6").Value = -60871.6
$ws.Range("H55").Value = 375.9375
$ws.Range("I55").Value = 238.63637
$ws.Range("K55").Value = 238.63637
$ws.Range("M55").Value = -65.63637
$ws.Range("H61").Value = 2006.4286
$ws.Range("I61").Value = 2115.3845
$ws.Range("J61").Value = 590
$ws.Range("K61").Value = 2115.3845
$ws.Range("L61").Value = 590
$ws.Range("M61").Value = -1913.3845
$ws.Range("N61").Value = -994
$ws.Range("H100").Value = 3215.8096
$ws.Range("I100").Value = 2636.6
$ws.Range("J100").Value = 14800
$ws.Range("K100").Value = 2636.6
$ws.Range("L100").Value = 14800
$ws.Range("M100").Value = -2095.6
$ws.Range("N100").Value = -15882
$ws.Range("H113").Value = 2006.4286
$ws.Range("I113").Value = 2115.3845
$ws.Range("J113").Value = 590
$ws.Range("K113").Value = 2115.3845
$ws.Range("L113").Value = 590
$ws.Range("M113").Value = 54.61549999999988
$ws.Range("N113").Value = -4930
$ws.Range("H132").Value = 8821.546
$ws.Range("I132").Value = 8378.286
$ws.Range("J132").Value = 9597.25
$ws.Range("K132").Value = 25134.858
$ws.Range("L132").Value = 28791.75
$ws.Range("M132").Value = -22604.858
$ws.Range("N132").Value = -33851.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 31686
$ws.Range("I54").Value = 30500
$ws.Range("J54").Value = 31982.5
$ws.Range("K54").Value = 30500
$ws.Range("L54").Value = 31982.5
$ws.Range("M54").Value = -29980
$ws.Range("N54").Value = -33022.5
$ws.Range("H74").Value = 4972.2
$ws.Range("J74").Value = 2560.5
$ws.Range("L74").Value = 2560.5
$ws.Range("N74").Value = -4432.5
$ws.Range("H77").Value = 4972.2
$ws.Range("J77").Value = 2560.5
$ws.Range("L77").Value = 7681.5
$ws.Range("N77").Value = -17041.5
